$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells P1 (14) and Q1 (15), matching the bold/centered/bordered
# style already used by the rest of row 1 (B1:O1).
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Updated B2:I25 block (columns J:O are untouched, still zero).
$data = @(
    @(24.44307359118704, 18.84945595138836, 10.02066305241336, 29.10526702065993, 63.03962927701833, 2.074527469938671, 3.089966362669392, 3.037267292020123),
    @(22.75403469952843, 17.53033396951019, 9.477274155190834, 27.02775996347093, 59.07048019432784, 2.084514762941586, 2.711583836633365, 2.666704839883651),
    @(21.65952378285824, 16.68690558818472, 9.126601439550928, 25.68821643124625, 56.50830951348823, 2.090765124095499, 2.473648940128331, 2.509508808433992),
    @(21.19839618377675, 16.35017261237931, 8.962248306382884, 25.1250067760214, 55.35697436933109, 2.093373239817404, 2.374520949779008, 2.613255113535315),
    @(21.12034009970148, 16.31276175098225, 8.913715193426158, 25.02928931015711, 55.07122712887607, 2.093843404499031, 2.357148475066212, 2.635588437125068),
    @(21.65169548959532, 16.73292013436325, 9.067926423508622, 25.67736305216642, 56.24481937174074, 2.090896024981767, 2.470075046960208, 2.524493452909088),
    @(23.870143067299, 18.46235242138411, 9.765799130928215, 28.39744764467337, 61.38995501941691, 2.078075537225312, 2.957617360055809, 2.908925381665187),
    @(27.77635781228673, 21.51752086636726, 11.09435833615976, 33.24354161310295, 70.89697344570428, 2.053570185398751, 3.877518438033591, 3.818003966742009),
    @(30.36896984993525, 23.60291155827736, 11.73454473692673, 35.62103610224544, 76.20990935265392, 2.036709757866656, 4.482358605421413, 4.449622519871673),
    @(31.36121202343197, 24.30366210427836, 9.970562816033523, 29.02521488677462, 69.78476296936714, 2.03570407727512, 4.790747882991368, 4.510909130059678),
    @(31.68585877720369, 24.44128593523203, 8.470649040072542, 23.0426619106675, 63.59460163662632, 2.037548408334782, 5.512775492837027, 4.456295562001463),
    @(31.52803904234584, 24.22003956114179, 7.009883965657012, 17.05336036085701, 56.73539812628097, 2.041798004963339, 6.450362237943956, 4.308155851271346),
    @(31.20074615630609, 23.91840789243268, 6.025607134096393, 12.91878041738623, 51.49898527829699, 2.045765259257151, 7.201276089694248, 4.16614545595082),
    @(31.01509405687257, 23.77999269105106, 5.778124340762592, 11.90505057293832, 49.98084492291337, 2.047349627796285, 7.37155626470092, 4.109105430901421),
    @(29.98938299365462, 23.00271204771266, 5.728103560500506, 11.57793924659197, 48.56665406277664, 2.053498711515744, 7.055752869198356, 3.877767251533557),
    @(29.36448092993919, 22.55793443121987, 6.19513540705, 13.59355848671587, 50.37157133274903, 2.056209579269872, 6.353811635908204, 3.772518732282649),
    @(29.0326828153557, 22.33313476896589, 7.282784649532969, 18.16992540045136, 55.3707338003543, 2.055887887705688, 5.333293727927019, 3.773529732633144),
    @(28.97345199471054, 22.4033883500889, 8.790375362995041, 24.53208797344972, 62.1806771869171, 2.052752511625638, 4.411514987705797, 3.879578676368524),
    @(29.69820850587216, 23.17917654727976, 11.41151704201456, 34.95363493126855, 74.17869503357267, 2.041462083656163, 4.314225926561374, 4.278525854025554),
    @(31.63458094084674, 24.75784781552005, 12.27927442675885, 38.09980501401512, 79.78419766314934, 2.027249542232134, 4.857352354950227, 4.805519619189777),
    @(32.85232260185697, 25.7093163221756, 12.77195792536501, 39.65220401541557, 83.08140148289172, 2.018208876805866, 5.187829086509634, 5.140100061090803),
    @(32.20893271165753, 25.16305894123705, 12.56410863771151, 38.83252030482961, 81.55397815487034, 2.022918405131028, 5.014325094539056, 4.96347957896903),
    @(29.68185813993916, 23.11238178875517, 11.66976099974378, 35.63223609375809, 75.24053863113274, 2.04088703992831, 4.351378196554124, 4.293550467567776),
    @(26.76275660477796, 20.79094544270774, 10.65831984704363, 31.97679880725834, 68.04404505137312, 2.060308407478158, 3.628166126820007, 3.571692568677724),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $col = $j + 2  # column B = 2
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}

# New columns P and Q for data rows 2:25, all zero (same as J:O).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 16).Value = 0  # P
    $ws.Cells.Item($r, 17).Value = 0  # Q
}
